# The author split the tail of a strikethrough sentence
#   " not cooking this up by myself. I learn this from this video:"
# by inserting the letter "t" after "I learn" (-> "I learnt"), turning the
# single run into three runs:
#   1) " not cooking this up by myself. I learn"   (untouched - keeps its
#      original w:rsidR/w:rsidRPr)
#   2) "t"                                          (new run)
#   3) " this from this video:"                     (rest of the original
#      text, re-run without rsid since it was re-written)
#
# Word's COM Range editing normally rebuilds (and drops the rsid of) any
# run a Find/Replace or InsertAfter/InsertBefore touches - even the
# untouched neighbour - unless the edit is a pure formatting round trip
# (apply then immediately revert a property) on text whose content does
# not change. We lean on that to keep run (1)'s rsid attributes intact
# while still slicing the paragraph into the three runs the diff expects.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("I learn this from this video:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target sentence to edit."
}

$start = $rng.Start
$end = $rng.End

# Step 1: split " not cooking this up by myself. I learn" | " this from
# this video:" right after "I learn" (start+7), turning the character
# that used to be the space into its own run, via a no-visible-effect
# Bold round trip so the engine is forced to materialize a run boundary
# there. Then overwrite that boundary run's text with "t" (dropping the
# space) - this mutates only the boundary run, leaving run (1) untouched
# so it keeps its original rsid attributes.
$boundary = $d.Range($start + 7, $start + 8)
$boundary.Font.Bold = 1
$boundary.Text = "t"
$boundary.Font.Bold = 0

# Step 2: restore the leading space onto the remainder of the sentence
# (now starting right after the new "t"), again wrapped in a harmless
# Bold round trip so it becomes its own run rather than re-merging with
# run (2).
$rest = $d.Range($start + 8, $end)
$rest.Font.Bold = 1
$rest.Text = " " + $rest.Text
$rest.Font.Bold = 0
